$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above the current row 4 (elasticity_gnrl_rate_occupancy_to_gdppc),
# shifting it and all rows below it down by one.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new climate-change-factor variable.
$ws.Range("A4").Value = "General"
$ws.Range("B4").Value = "climate_change_factor_gnrl_hydropower_availability"
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 0.5
$ws.Range("J4:AS4").Value = 1
